# Generate Report for Handoff
#
# The localization status report was regenerated: the handed-off file's
# GUID-based name changed (and the content hash of the generated xliff
# files derived from it), and the associated handoff / xliff-generate
# timestamps moved a few seconds later. Apply the same text
# substitutions the report generator produced, across the Overview,
# zh-cn and de-de sheets, and refresh the hyperlink display text to
# match the new file name (the hyperlink target itself is left as-is).

$wb = $excel.ActiveWorkbook

$newGuid = "abf1ff20-febe-4748-bfb2-66b3e1478b0d"
$newHash = "54ad1d06a6b0dd538de3b4949ff8986929e56b4b"

$hyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4cc3b25db8e7d1181f5fa44655a34abf4059cfb7/e2e/c0e76c29-e479-4ae7-b530-ab601efbd79b.md"

# Hyperlink colour/underline used throughout this workbook, re-applied
# after recreating a hyperlink so its look matches the rest of the sheet.
$hyperlinkColor = 15570276  # OLE (BGR) form of RGB 0x6495ED

function Update-HandoffHyperlink($ws, $cellAddr, $displayText) {
    $rng = $ws.Range($cellAddr)
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($rng, $hyperlinkUrl, "", "", $displayText)
    $rng.Font.Color = $hyperlinkColor
    $rng.Font.Underline = 2
}

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-24 20:59:54"
Update-HandoffHyperlink $wsOverview "B2" "e2e\$newGuid.md"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-24 20:59:49"
Update-HandoffHyperlink $wsZhCn "A2" "$newGuid.md"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-24 20:59:54"
Update-HandoffHyperlink $wsDeDe "A2" "$newGuid.md"
